$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The template already has one "Шмурдяк" category block in rows 11-13
# (Изабелла/Шардоне/Гранатовый браслет). Add a marketing tag by duplicating
# that block into new rows 14-16 so the tag/category repeats further down
# the template.

# Copy values from the source block into the new rows
$ws.Range("A11:E13").Copy()
$ws.Range("A14").Select()
$ws.Paste()
$excel.CutCopyMode = 0

# Make sure the new rows carry over the exact same cell formatting (style)
# as the source rows
$ws.Range("A11:E13").Copy()
$ws.Range("A14:E16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match row heights to the source rows as well
$ws.Rows.Item(14).RowHeight = $ws.Rows.Item(11).RowHeight
$ws.Rows.Item(15).RowHeight = $ws.Rows.Item(12).RowHeight
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(13).RowHeight

# Update the active selection to the newly added block, matching the
# workbook's saved view state
$ws.Range("A14:E16").Select()
